$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.716.04"
$ws.Range("E2").Value = "  +0.82%  "

$ws.Range("D3").Value = "1.771.41"
$ws.Range("E3").Value = "  +1.34%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.92"
$ws.Range("E5").Value = "  +0.73%  "

$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4589"
$ws.Range("E7").Value = "  +3.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3579"
$ws.Range("E8").Value = "  -0.57%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07479"
$ws.Range("E9").Value = "  -0.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.82"
$ws.Range("E10").Value = "  -0.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.099"
$ws.Range("E11").Value = "  +0.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.01%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.78"
$ws.Range("E13").Value = "  +0.98%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.031"
$ws.Range("E14").Value = "  +0.20%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.209"
$ws.Range("E15").Value = "  +1.23%  "

$ws.Range("D16").Value = "1.773.29"
$ws.Range("E16").Value = "  +1.20%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.57"
$ws.Range("E17").Value = "  +0.63%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001058"
$ws.Range("E18").Value = "  -0.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06415"
$ws.Range("E19").Value = "  +0.21%  "

$ws.Range("E20").Value = "  -0.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.03"
$ws.Range("E21").Value = "  +1.56%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.789"
$ws.Range("E22").Value = "  -0.98%  "

$ws.Range("D23").Value = "27.782.51"
$ws.Range("E23").Value = "  +0.87%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.28"
$ws.Range("E24").Value = "  +1.16%  "

$ws.Range("E25").Value = "  -0.84%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.17"
$ws.Range("E26").Value = "  +1.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.22"
$ws.Range("E27").Value = "  -1.10%  "

$ws.Range("D28").Value = "1.977.12"
$ws.Range("E28").Value = "  +1.20%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.171"
$ws.Range("E29").Value = "  +4.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.71"
$ws.Range("E30").Value = "  +0.64%  "

$ws.Range("E31").Value = "  +0.49%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09214"
$ws.Range("E32").Value = "  +2.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.672"
$ws.Range("E33").Value = "  +0.39%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.516"
$ws.Range("E34").Value = "  -0.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.82"
$ws.Range("E35").Value = "  -1.33%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02292"
$ws.Range("E36").Value = "  -0.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06171"
$ws.Range("E37").Value = "  +2.62%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2085"
$ws.Range("E38").Value = "  +0.15%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6307"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.938"
$ws.Range("E40").Value = "  -0.15%  "

$ws.Range("E41").Value = "  -1.95%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.391"
$ws.Range("E42").Value = "  +0.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.765"
$ws.Range("E43").Value = "  +0.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.19"
$ws.Range("E44").Value = "  +0.61%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.741"
$ws.Range("E45").Value = "  +0.80%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5885"
$ws.Range("E46").Value = "  +0.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.36"
$ws.Range("E47").Value = "  +0.28%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.946"
$ws.Range("E48").Value = "  -0.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06926"
$ws.Range("E49").Value = "  +1.12%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.134"
$ws.Range("E50").Value = "  -0.86%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.09"
$ws.Range("E51").Value = "  +0.37%  "
